$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.182.13'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '3.807.54'
$ws.Range("E3").Value = '  +1.71%  '

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E4").Value = '  +0.06%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '602.62'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.01%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '163.90'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  -2.86%  '

$ws.Range("D7").Value = '3.805.58'
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("E10").Value = '  +2.37%  '

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.31'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("E12").Value = '  -0.47%  '

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '37.37'
$cell.Style = $origStyle
$ws.Range("E13").Value = '  -2.28%  '

$ws.Range("E14").Value = '  -0.88%  '

$ws.Range("D15").Value = '4.445.67'
$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("D16").Value = '3.802.61'
$ws.Range("E16").Value = '  +1.39%  '

$ws.Range("D17").Value = '69.287.38'
$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("E18").Value = '  +1.96%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.37'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +1.45%  '

$ws.Range("E20").Value = '  -0.33%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.28'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +4.07%  '

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '490.91'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  -0.66%  '

$ws.Range("E23").Value = '  -0.51%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000152'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  -2.09%  '

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '84.72'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  -0.67%  '

$ws.Range("E26").Value = '  -2.95%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.26'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  -0.77%  '

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.08'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  -3.25%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("E30").Value = '  -0.01%  '

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.08'
$cell.Style = $origStyle
$ws.Range("E31").Value = '  +1.51%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.40'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  -5.07%  '

$ws.Range("D33").Value = '3.953.60'
$ws.Range("E33").Value = '  +1.66%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '32.04'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  +0.68%  '

$ws.Range("D35").Value = '3.755.41'
$ws.Range("E35").Value = '  +2.14%  '

$ws.Range("E36").Value = '  -1.70%  '

$ws.Range("E37").Value = '  +5.89%  '

$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("E39").Value = '  +1.18%  '

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  +0.11%  '

$ws.Range("E41").Value = '  +0.12%  '

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.04'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +1.04%  '

$ws.Range("E43").Value = '  +0.87%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '48.51'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  -0.76%  '

$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '422.86'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  -3.22%  '

$ws.Range("E46").Value = '  +0.00%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.40'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  -0.65%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.830.12'
$ws.Range("E48").Value = '  +2.05%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '142.05'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +0.33%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '39.70'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  -1.78%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.28'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +4.81%  '
